$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Gran Canaria overtook A Coruña in "Casos totales" -> rows swap places
# while sorted descending by total cases (col B).
$ws.Range("A19").Value = "Gran Canaria"
$ws.Range("B19").Value = 878
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 818
$ws.Range("E19").Value = 27

$ws.Range("A20").Value = "A Coruña"
$ws.Range("B20").Value = 824
$ws.Range("C20").Value = 47
$ws.Range("D20").Value = 789
$ws.Range("E20").Value = 23

# Death-count ("Muertes") corrections for several other provinces.
$ws.Range("E31").Value = 27
$ws.Range("E55").Value = 27
$ws.Range("E57").Value = 27
$ws.Range("E58").Value = 27
$ws.Range("E62").Value = 27
$ws.Range("E63").Value = 27

# Timestamp footer update.
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 22:42"
